# Remove the DNP resistor (R1) and the mounting-hole (MH1-MH4) rows from
# the BOM on the "DAC_bom_qty_10" worksheet.
#
# In the original sheet:
#   Row 12 = Qty 4, Designator "MH1, MH2, MH3, MH4", Comment "MOUNT_HOLE"
#   Row 13 = Qty 1, Designator "R1", Comment "Resistor", Value "DNP"
# Both rows are deleted entirely, and every row below shifts up by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAC_bom_qty_10")
$ws.Activate()

# Delete the two BOM rows (MOUNT_HOLE and DNP Resistor) - remaining rows
# automatically shift up to fill the gap.
$ws.Rows("12:13").Delete()

# Leave the view focused on the first row that follows the deleted rows,
# matching where the edit took place.
$ws.Rows("12:12").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
